# Actualización desde MV -datos-
# Adds the daily "Tasa de Política Monetaria" readings for 20-09-2021
# (completing the row that was previously partial) plus five new trading
# days: 21, 22, 23, 24 and 27 September 2021.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 188 (20-09-2021) — fill in the remaining series values ---
$ws.Range("D188").Value = -0.1
$ws.Range("E188").Value = -0.5
$ws.Range("H188").Value = 1.75
$ws.Range("M188").Value = 1.13
$ws.Range("O188").Value = 38
$ws.Range("P188").Value = 5.25
$ws.Range("Q188").Value = 1.75

# --- Row 189 (21-09-2021) ---
$ws.Range("A189").Value = "21-09-2021"
$ws.Range("B189").Value = 0.25
$ws.Range("C189").Value = 0.1
$ws.Range("D189").Value = -0.1
$ws.Range("E189").Value = -0.5
$ws.Range("F189").Value = 0.75
$ws.Range("H189").Value = 1.75
$ws.Range("I189").Value = 0.1
$ws.Range("J189").Value = 0.75
$ws.Range("K189").Value = 6.75
$ws.Range("L189").Value = 0.5
$ws.Range("M189").Value = 1.13
$ws.Range("N189").Value = 19
$ws.Range("O189").Value = 38
$ws.Range("P189").Value = 5.25
$ws.Range("Q189").Value = 1.75
$ws.Range("R189").Value = 4.5
$ws.Range("S189").Value = 1

# --- Row 190 (22-09-2021) ---
$ws.Range("A190").Value = "22-09-2021"
$ws.Range("B190").Value = 0.25
$ws.Range("C190").Value = 0.1
$ws.Range("D190").Value = -0.1
$ws.Range("E190").Value = -0.5
$ws.Range("F190").Value = 0.75
$ws.Range("G190").Value = 4.35
$ws.Range("H190").Value = 1.75
$ws.Range("I190").Value = 0.1
$ws.Range("J190").Value = 0.75
$ws.Range("K190").Value = 6.75
$ws.Range("L190").Value = 0.5
$ws.Range("M190").Value = 1.13
$ws.Range("N190").Value = 19
$ws.Range("O190").Value = 38
$ws.Range("P190").Value = 6.25
$ws.Range("Q190").Value = 1.75
$ws.Range("R190").Value = 4.5
$ws.Range("S190").Value = 1

# --- Row 191 (23-09-2021) ---
$ws.Range("A191").Value = "23-09-2021"
$ws.Range("B191").Value = 0.25
$ws.Range("C191").Value = 0.1
$ws.Range("D191").Value = -0.1
$ws.Range("E191").Value = -0.5
$ws.Range("G191").Value = 4.35
$ws.Range("H191").Value = 1.75
$ws.Range("I191").Value = 0.1
$ws.Range("J191").Value = 0.75
$ws.Range("K191").Value = 6.75
$ws.Range("L191").Value = 0.5
$ws.Range("M191").Value = 1.13
$ws.Range("N191").Value = 18
$ws.Range("O191").Value = 38
$ws.Range("P191").Value = 6.25
$ws.Range("Q191").Value = 1.75
$ws.Range("R191").Value = 4.5
$ws.Range("S191").Value = 1

# --- Row 192 (24-09-2021) ---
$ws.Range("A192").Value = "24-09-2021"
$ws.Range("B192").Value = 0.25
$ws.Range("C192").Value = 0.1
$ws.Range("D192").Value = -0.1
$ws.Range("E192").Value = -0.5
$ws.Range("F192").Value = 0.75
$ws.Range("G192").Value = 4.35
$ws.Range("H192").Value = 1.75
$ws.Range("I192").Value = 0.1
$ws.Range("J192").Value = 0.75
$ws.Range("K192").Value = 6.75
$ws.Range("L192").Value = 0.5
$ws.Range("M192").Value = 1.13
$ws.Range("N192").Value = 18
$ws.Range("O192").Value = 38
$ws.Range("P192").Value = 6.25
$ws.Range("Q192").Value = 1.75
$ws.Range("R192").Value = 4.5
$ws.Range("S192").Value = 1

# --- Row 193 (27-09-2021) — partial row, only through column N ---
$ws.Range("A193").Value = "27-09-2021"
$ws.Range("C193").Value = 0.1
$ws.Range("F193").Value = 0.75
$ws.Range("G193").Value = 4.35
$ws.Range("I193").Value = 0.1
$ws.Range("J193").Value = 0.75
$ws.Range("K193").Value = 6.75
$ws.Range("L193").Value = 0.5
$ws.Range("N193").Value = 18
